$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 83, shifting existing rows 83..190 down to 84..191
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new data record
$ws.Cells.Item(83, 1).Value  = 7
$ws.Cells.Item(83, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(83, 3).Value  = "Ñuble"
$ws.Cells.Item(83, 4).Value  = 44601
$ws.Cells.Item(83, 5).Value  = 16
$ws.Cells.Item(83, 6).Value  = 100112032
$ws.Cells.Item(83, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(83, 8).Value  = "Sin especificar"
$ws.Cells.Item(83, 9).Value  = "Primera"
$ws.Cells.Item(83, 10).Value = 120
$ws.Cells.Item(83, 11).Value = 7500
$ws.Cells.Item(83, 12).Value = 8000
$ws.Cells.Item(83, 13).Value = 7750
$ws.Cells.Item(83, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(83, 15).Value = "Región del Maule"
$ws.Cells.Item(83, 16).Value = 155
$ws.Cells.Item(83, 17).Value = 50
$ws.Cells.Item(83, 18).Value = "Hortaliza"
